$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "alamat" column: move no_telepon/email one column to the left
# (into C1/D1), then clear the now-trailing E1 header and the old C2 note.
$ws.Range("C1").Value2 = $ws.Range("D1").Value2
$ws.Range("D1").Value2 = $ws.Range("E1").Value2
$ws.Range("E1").ClearContents()
$ws.Range("C2").ClearContents()

# Selection ends up at C6 in the saved file.
$ws.Range("C6").Select() | Out-Null
